$d = $word.ActiveDocument

# --- Step 1: mint numId=2 (bullet abstractNum) and numId=3 (decimal/lowerLetter abstractNum) ---
# Apply on isolated paragraphs far apart so the runtime mints two distinct list definitions,
# matching the target's numId=2 (bullet) / numId=3 (decimal) layout.
$mint1 = $d.Paragraphs(1)
$mint1.Range.ListFormat.ApplyBulletDefault()
$mint2 = $d.Paragraphs($d.Paragraphs.Count)
$mint2.Range.ListFormat.ApplyNumberDefault()

# --- Step 2: replace the entire body content with the target OOXML ---
$xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Home</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>“Only after the last tree has been cut down, the last river has been poisoned, and the last fish has been caught, only then will people realize that money cannot be eaten.”</w:t></w:r><w:r><w:br/><w:t>“Look deep into nature and then you will understand everything better”</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The Earth is the only planet we humans can call home, and as far as we know this is the only planet that can sustain life. Every day is a constant struggle for the planet to keep up with our modern lifestyle. We humans who have been utilizing and misusing nature’s resources and benefitting from it must also learn to protect it and nurture it. Our advocacy leans towards eco-friendly projects and acts of reducing wastes and other resources that might pose a threat to the environment and human health. We promote recycling and proper garbage disposal in accordance to RA 6969 </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:bCs/><w:i/></w:rPr><w:t xml:space="preserve">Republic Act 6969 Toxic Substances, Hazardous and Nuclear Waste Control Act </w:t></w:r><w:r><w:rPr><w:bCs/><w:i/></w:rPr><w:t>of</w:t></w:r><w:r><w:rPr><w:bCs/><w:i/></w:rPr><w:t xml:space="preserve"> 1990</w:t></w:r><w:r><w:rPr><w:bCs/><w:i/></w:rPr><w:t xml:space="preserve">). </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="720"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Laws that Support Our Cause</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:bCs/></w:rPr><w:t>“REPUBLIC ACT 6969 TOXIC SUBSTANCES, HAZARDOUS AND NUCLEAR WASTE CONTROL ACT OF 1990</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>The law aims to regulate restrict or prohibit the importation, manufacture, processing, sale, distribution, use and disposal of chemical substances and mixtures the present unreasonable risk to human health. It territor</w:t></w:r><w:r><w:t>ial l</w:t></w:r><w:r><w:t xml:space="preserve">ikewise prohibits the entry, even in transit, of hazardous and nuclear wastes and their disposal into the Philippine </w:t></w:r><w:r><w:t>l</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>imits for whatever purpose; and to provide advancement and facilitate research and studies on toxic chemicals.”</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t>REPUBLIC ACT 9003 ECOLOGICAL SOLID WASTE MANAGEMENT ACT OF 2000</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>In partnership with stakeholders, the law aims to adopt a systematic, comprehensive and ecological solid waste management program that shall ensure the protection of public health and environment. The law ensures proper segregation, collection, storage, treatment and disposal of solid waste through the formulation and adaptation of best eco-waste products</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">PRESIDENTIAL DECREE 1586 ENVIRONMENTAL IMPACT STATEMENT (EIS) STATEMENT OF 1978 </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>The Environment Impact Assessment System was formally established in 1978 with the enactment of Presidential Decree no. 1586 to facilitate the attainment and maintenance of rational and orderly balance between socio-economic development and environmental protection. EIA is a planning and management tool that will help government, decision makers, the proponents and the affected community address the negative consequences or risks on the environment. The process assures implementation of environment-friendly projects.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Vision-Mission</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">We envision the barangay of Bakakeng Norte to be a clean and green community, that upholds the values and practices for protecting and preserving the environment that respects the mountains ecosystem. We see to it that no resources go to waste, and have the community participate and help with our cause. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Be aware</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:t xml:space="preserve">According to the interview conducted from the Barangay officials in Bakakeng Norte, one of the aforementioned problems concerning the environment is that some households do not segregate their garbage, this can be a handful especially for the garbage collectors who have to sort out the garbage beforehand. The garbage disposal schedule is not being followed by some residents and at in some days’ garbage collection tend to be late.  Littering is also another issue that was mentioned by the interviewed officials. On a minor scale one of the issues that was raised by one of the officials is the burning of garbage that occurs occasionally on some houses, this can pose a dangerous threat to resident’s health if inhaled and should be stopped and condoned. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Help out in your own way</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>In your household:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Segregate your waste</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Have your home practice proper waste disposal by segregating your trash into non-biodegradable</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>(di-nabubulok)</w:t></w:r><w:r><w:t xml:space="preserve">, biodegradable trash </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>(nabubulok)</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>Your biodegradable trash can serve the purpose of being a compost for your garden and will serve as a natural fertilizer. While your non-biodegradable waste can be segregated further into different categories namely:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Paper</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Cardboard</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Glass</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Plastic Bottles</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Plant a tree or plants that can absorb carbon monoxide than ordinary household plants</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Do not waste water, and make sure to avoid leaky faucets</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Save power by turning off your lights during the day or when not in use. This can help save you from spending too much on electricity plus you’re able to contribute to our cause.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Follow the barangay’s garbage disposal schedule. This way dogs won’t end up scattering the trash from the trash pile and it won’t create a pungent smell that might inconvenience the community.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Do not burn your garbage, this is not the proper way to dispose of your garbage. In fact, you are creating toxic fumes that will be dangerous to both your health and the people nearby who is inhaling the smoke coming from the burnt waste.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">When shopping for grocery or household items, try to be minimal about using plastic “sando” bags, and use an ecofriendly bag instead. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Avoid using products that contain UFCs, these contains chemicals that is affecting the Earth’s ozone layer.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>If you are in position of a vehicle, please have a biweekly emission testing to minimize the smog created by your vehicle.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>We encourage the community to sought renewable power sources to avoid wasting resources.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Follow the 4Rs</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Reduce</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Reuse</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Recycle</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Repeat</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="14400"/></w:pPr><w:r><w:br w:type="page"/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>Resources:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>https://ecac.emb.gov.ph/?page_id=43</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>http://119.92.161.2/embgovph/ecac/LGUCorner/MajorEnvironmentalLaws/PresidentialDecree1586EnvironmentalImpactStat.aspx</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
"@
$d.Content.InsertXML($xml)

Write-Host "Paragraphs after edit: $($d.Paragraphs.Count)"
